# Applies the OOXML diff: splits the tab-run that precedes "DATE:"/"Date:"
# in the two booking-letter paragraphs into several tab-only runs plus a
# couple of extra spacing runs, matching the commit's "forwading templets"
# realignment.

$d = $word.ActiveDocument

# ---- Paragraph 1: "Ref :NO.MSL/${msl_num}/${year} ... DATE: ... " ----
$para1New = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000E0CB2" w:rsidRDefault="008032F7" w:rsidP="000E0CB2">' +
    '<w:r><w:t xml:space="preserve">Ref </w:t></w:r>' +
    '<w:r w:rsidR="00C51341"><w:t>:</w:t></w:r>' +
    '<w:r><w:t>NO</w:t></w:r>' +
    '<w:r w:rsidR="007F23FD"><w:t>.</w:t></w:r>' +
    '<w:r w:rsidR="00C47E32"><w:t>MSL/</w:t></w:r>' +
    '<w:r w:rsidR="000E7293"><w:t>${</w:t></w:r>' +
    '<w:r w:rsidR="00491098"><w:t>msl_num</w:t></w:r>' +
    '<w:r w:rsidR="00B16212"><w:t>}</w:t></w:r>' +
    '<w:r w:rsidR="00493084"><w:t>/</w:t></w:r>' +
    '<w:r w:rsidR="00B16212"><w:t>${year}</w:t></w:r>' +
    '<w:r w:rsidR="00397550"><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidR="00397550"><w:tab/></w:r>' +
    '<w:r w:rsidR="00397550"><w:tab/></w:r>' +
    '<w:r w:rsidR="00397550"><w:tab/></w:r>' +
    '<w:r><w:tab/><w:t xml:space="preserve">          </w:t></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:tab/><w:t xml:space="preserve">   </w:t></w:r>' +
    '<w:r><w:t>DATE:        .${month}.${year}</w:t></w:r>' +
    '</w:p>'

# ---- Paragraph 2: "Ref: MSL/${msl_num}/${year} ... Date: ... " ----
$para2New = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000E0CB2" w:rsidRDefault="00EC0BE5" w:rsidP="00397550">' +
    '<w:pPr><w:jc w:val="both"/></w:pPr>' +
    '<w:r><w:t>Ref: MSL/</w:t></w:r>' +
    '<w:r w:rsidR="000E7293"><w:t>${</w:t></w:r>' +
    '<w:r w:rsidR="00491098"><w:t>msl_num</w:t></w:r>' +
    '<w:r w:rsidR="00B16212"><w:t>}</w:t></w:r>' +
    '<w:r w:rsidR="004277FB"><w:t>/</w:t></w:r>' +
    '<w:r w:rsidR="00B16212"><w:t>${year}</w:t></w:r>' +
    '<w:r w:rsidR="00397550" w:rsidRPr="00397550"><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidR="00397550"><w:tab/></w:r>' +
    '<w:r w:rsidR="00397550"><w:tab/></w:r>' +
    '<w:r w:rsidR="00397550"><w:tab/></w:r>' +
    '<w:r w:rsidR="00397550"><w:tab/></w:r>' +
    '<w:r w:rsidR="00397550"><w:tab/></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:tab/><w:t xml:space="preserve">      </w:t></w:r>' +
    '<w:r w:rsidR="00397550" w:rsidRPr="007B5C2F"><w:t>Date:</w:t></w:r>' +
    '<w:r w:rsidR="00397550"><w:t xml:space="preserve">        </w:t></w:r>' +
    '<w:r w:rsidR="00397550" w:rsidRPr="007B5C2F"><w:t>.</w:t></w:r>' +
    '<w:r w:rsidR="00397550"><w:t>${month}.${year}</w:t></w:r>' +
    '</w:p>'

$count = $d.Paragraphs.Count
$done1 = $false
$done2 = $false

for ($i = 1; $i -le $count; $i++) {
    if ($done1 -and $done2) { break }

    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ((-not $done1) -and $t.Contains('DATE:') -and $t.Contains('MSL/')) {
        $r = $p.Range
        $r.Text = ""
        [void]$r.InsertXML($para1New)
        $done1 = $true
        continue
    }

    if ((-not $done2) -and $t.Contains('Date:') -and $t.Contains('MSL/')) {
        $r = $p.Range
        $r.Text = ""
        [void]$r.InsertXML($para2New)
        $done2 = $true
        continue
    }
}

Write-Host "Paragraph 1 updated:" $done1
Write-Host "Paragraph 2 updated:" $done2
